$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update of price/quality data for rows 5-20 (D, I, J, K, L, M, P columns).
# Values below represent the new "after" state for each row, derived from the commit diff.

$ws.Range("D5").Value = 44245
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 12000
$ws.Range("P5").Value = 667
$ws.Range("D6").Value = 44245
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 10000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 10000
$ws.Range("P6").Value = 556
$ws.Range("D7").Value = 44229
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 15000
$ws.Range("P7").Value = 833
$ws.Range("D8").Value = 44383
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 16000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 16000
$ws.Range("P8").Value = 889
$ws.Range("D9").Value = 44383
$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 12000
$ws.Range("P9").Value = 667
$ws.Range("D10").Value = 44238
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 12000
$ws.Range("P10").Value = 667
$ws.Range("D11").Value = 44238
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 10000
$ws.Range("P11").Value = 556
$ws.Range("D12").Value = 44238
$ws.Range("I12").Value = "Tercera"
$ws.Range("J12").Value = 50
$ws.Range("K12").Value = 8000
$ws.Range("L12").Value = 8000
$ws.Range("M12").Value = 8000
$ws.Range("P12").Value = 444
$ws.Range("D13").Value = 44235
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 13000
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = 13000
$ws.Range("P13").Value = 722
$ws.Range("D14").Value = 44235
$ws.Range("I14").Value = "Segunda"
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 11000
$ws.Range("L14").Value = 11000
$ws.Range("M14").Value = 11000
$ws.Range("P14").Value = 611
$ws.Range("D15").Value = 44235
$ws.Range("I15").Value = "Tercera"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 9000
$ws.Range("L15").Value = 9000
$ws.Range("M15").Value = 9000
$ws.Range("P15").Value = 500
$ws.Range("D16").Value = 44391
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 15000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 15000
$ws.Range("P16").Value = 833
$ws.Range("D17").Value = 44249
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 400
$ws.Range("K17").Value = 12000
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = 12000
$ws.Range("P17").Value = 667
$ws.Range("D18").Value = 44249
$ws.Range("I18").Value = "Segunda"
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = 10000
$ws.Range("P18").Value = 556
$ws.Range("D19").Value = 44396
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = 15000
$ws.Range("L19").Value = 15000
$ws.Range("M19").Value = 15000
$ws.Range("P19").Value = 833
$ws.Range("D20").Value = 44396
$ws.Range("I20").Value = "Segunda"
$ws.Range("J20").Value = 150
$ws.Range("K20").Value = 12000
$ws.Range("L20").Value = 12000
$ws.Range("M20").Value = 12000
$ws.Range("P20").Value = 667

Write-Output "Done updating rows 5-20"
